$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("E2").Value = "2024.03.16 09:00-03.17 17:00"
$ws.Range("F2").Value = 2789
$ws.Range("E3").Value = "2024.03.16 12:30-03.16 16:30"
$ws.Range("E4").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("E5").Value = "2024.03.16 09:00-03.17 17:00"
$ws.Range("F5").Value = 1551
$ws.Range("E6").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("F6").Value = 1150
$ws.Range("E7").Value = "2024.03.16 09:30-03.16 16:00"
$ws.Range("E8").Value = "2024.03.16 09:00-03.16 17:00"
$ws.Range("F8").Value = 547
$ws.Range("E9").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("E10").Value = "2024.03.17 12:30-03.17 16:30"
$ws.Range("E11").Value = "2024.03.23 10:00-03.24 17:00"
$ws.Range("F11").Value = 9427
$ws.Range("E12").Value = "2024.03.23 10:00-03.23 17:00"
$ws.Range("F12").Value = 404
$ws.Range("E13").Value = "2024.03.23 09:30-03.23 17:00"
$ws.Range("F13").Value = 2508
$ws.Range("E14").Value = "2024.03.23 12:00-03.23 21:00"
$ws.Range("F14").Value = 11
$ws.Range("E15").Value = "2024.03.24 10:00-03.24 17:00"
$ws.Range("F15").Value = 264
$ws.Range("E16").Value = "2024.03.24 12:00-03.24 16:00"
$ws.Range("F16").Value = 183
$ws.Range("E17").Value = "2024.03.24 09:30-03.24 17:00"
$ws.Range("E18").Value = "2024.03.30 10:00-03.31 17:30"
$ws.Range("F18").Value = 656
$ws.Range("E19").Value = "2024.03.30 10:00-03.31 18:00"
$ws.Range("E20").Value = "2024.03.30 10:00-03.31 17:00"
$ws.Range("F20").Value = 1192
$ws.Range("E21").Value = "2024.03.30 10:00-03.31 17:00"
$ws.Range("E22").Value = "2024.04.04 09:30-04.05 16:30"
$ws.Range("F22").Value = 2933
$ws.Range("E23").Value = "2024.04.04 10:00-04.05 17:00"
$ws.Range("F23").Value = 2222
$ws.Range("E24").Value = "2024.04.04 10:00-04.05 17:00"
$ws.Range("E25").Value = "2024.04.04 10:00-04.05 17:00"
$ws.Range("E26").Value = "2024.04.05 09:30-04.05 16:30"
$ws.Range("F26").Value = 1934
$ws.Range("E27").Value = "2024.04.05 09:30-04.05 16:30"
$ws.Range("E28").Value = "2024.04.05 11:00-04.06 17:00"
$ws.Range("F28").Value = 1553
$ws.Range("E29").Value = "2024.04.13 09:00-04.13 18:00"
$ws.Range("E30").Value = "2024.04.13 13:00-04.13 19:00"
$ws.Range("F30").Value = 10
$ws.Range("E31").Value = "2024.04.13 10:00-04.13 17:00"
$ws.Range("E32").Value = "2024.04.13 09:00-04.13 18:00"
$ws.Range("F32").Value = 218
$ws.Range("E33").Value = "2024.04.20 10:00-04.21 17:00"
$ws.Range("E34").Value = "2024.04.20 09:00-04.20 22:00"
$ws.Range("F34").Value = 338
$ws.Range("E35").Value = "2024.04.20 14:00-04.20 18:00"
$ws.Range("E36").Value = "2024.04.20 10:00-04.21 18:00"
$ws.Range("E37").Value = "2024.04.20 10:00-04.20 17:00"
$ws.Range("F37").Value = 502
$ws.Range("E38").Value = "2024.04.30 10:00-05.03 16:00"
$ws.Range("E39").Value = "2024.05.01 10:00-05.02 17:00"
$ws.Range("F39").Value = 99
$ws.Range("E40").Value = "2024.05.01 10:30-05.02 17:00"
$ws.Range("F40").Value = 1266
$ws.Range("E41").Value = "2024.05.02 10:00-05.02 17:00"
$ws.Range("F41").Value = 98
$ws.Range("E42").Value = "2024.05.02 10:00-05.03 17:00"
$ws.Range("F42").Value = 1448
$ws.Range("E43").Value = "2024.05.02 10:00-05.02 17:00"
$ws.Range("F43").Value = 18
$ws.Range("E44").Value = "2024.05.18 10:00-05.18 21:00"
$ws.Range("F44").Value = 328
$ws.Range("E45").Value = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F45").Value = 18
$ws.Range("E46").Value = "2024.06.08 10:00-06.10 17:00"
$ws.Range("F46").Value = 199
$ws.Range("E47").Value = "2024.06.09 10:00-06.09 23:00"
$ws.Range("F47").Value = 706
$ws.Range("E48").Value = "2024.07.20 13:00-07.20 17:00"
$ws.Range("E49").Value = "2024.07.20 13:00-07.20 19:00"
$ws.Range("F49").Value = 310

$ws = $wb.Worksheets.Item(2)
$ws.Range("E2").Value = "2024.03.16 19:00-03.16 21:00"
$ws.Range("E3").Value = "2024.03.23 19:30-03.23 21:00"
$ws.Range("E4").Value = "2024.03.23 19:30-03.23 21:00"
$ws.Range("E5").Value = "2024.04.21 19:30-04.21 21:00"
$ws.Range("E6").Value = "2024.04.21 19:30-04.21 21:20"
$ws.Range("E7").Value = "2024.04.27 19:30-04.27 21:30"
$ws.Range("E8").Value = "2024.05.01 19:30-05.01 21:00"
$ws.Range("E9").Value = "2024.05.12 19:30-05.12 21:30"
$ws.Range("E10").Value = "2024.11.01 19:30-11.03 17:00"

$ws = $wb.Worksheets.Item(4)
$ws.Range("E2").Value = "2024.03.16 09:00-03.17 17:00"
$ws.Range("F2").Value = 2789
$ws.Range("E3").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("E4").Value = "2024.03.16 09:00-03.17 17:00"
$ws.Range("F4").Value = 1551
$ws.Range("E5").Value = "2024.03.16 19:00-03.16 21:00"
$ws.Range("E6").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("F6").Value = 1150
$ws.Range("E7").Value = "2024.03.16 09:00-03.16 17:00"
$ws.Range("F7").Value = 547
$ws.Range("E8").Value = "2024.03.16 10:00-03.16 17:00"
$ws.Range("E9").Value = "2024.03.23 10:00-03.24 17:00"
$ws.Range("F9").Value = 9427
$ws.Range("E10").Value = "2024.03.23 10:00-03.23 17:00"
$ws.Range("F10").Value = 404
$ws.Range("E11").Value = "2024.03.23 19:30-03.23 21:00"
$ws.Range("E12").Value = "2024.03.23 12:00-03.23 21:00"
$ws.Range("F12").Value = 11
$ws.Range("E13").Value = "2024.03.23 19:30-03.23 21:00"
$ws.Range("E14").Value = "2024.03.24 10:00-03.24 17:00"
$ws.Range("F14").Value = 264
$ws.Range("E15").Value = "2024.03.24 12:00-03.24 16:00"
$ws.Range("F15").Value = 183
$ws.Range("E16").Value = "2024.03.30 10:00-03.31 17:30"
$ws.Range("F16").Value = 656
$ws.Range("E17").Value = "2024.03.30 10:00-03.31 17:00"
$ws.Range("F17").Value = 1192
$ws.Range("E18").Value = "2024.03.30 10:00-03.31 17:00"
$ws.Range("E19").Value = "2024.04.04 09:30-04.05 16:30"
$ws.Range("F19").Value = 2933
$ws.Range("E20").Value = "2024.04.04 10:00-04.05 17:00"
$ws.Range("F20").Value = 2222
$ws.Range("E21").Value = "2024.04.04 10:00-04.05 17:00"
$ws.Range("E22").Value = "2024.04.05 09:30-04.05 16:30"
$ws.Range("E23").Value = "2024.04.05 11:00-04.06 17:00"
$ws.Range("F23").Value = 1553
$ws.Range("E24").Value = "2024.04.13 09:00-04.13 18:00"
$ws.Range("E25").Value = "2024.04.13 13:00-04.13 19:00"
$ws.Range("F25").Value = 10
$ws.Range("E26").Value = "2024.04.13 10:00-04.13 17:00"
$ws.Range("E27").Value = "2024.04.13 09:00-04.13 18:00"
$ws.Range("F27").Value = 218
$ws.Range("E28").Value = "2024.04.20 10:00-04.21 17:00"
$ws.Range("E29").Value = "2024.04.20 09:00-04.20 22:00"
$ws.Range("F29").Value = 338
$ws.Range("E30").Value = "2024.04.20 14:00-04.20 18:00"
$ws.Range("E31").Value = "2024.04.20 10:00-04.21 18:00"
$ws.Range("E32").Value = "2024.04.20 10:00-04.20 17:00"
$ws.Range("F32").Value = 502
$ws.Range("E33").Value = "2024.04.21 19:30-04.21 21:00"
$ws.Range("E34").Value = "2024.04.21 19:30-04.21 21:20"
$ws.Range("E35").Value = "2024.04.27 19:30-04.27 21:30"
$ws.Range("E36").Value = "2024.04.30 10:00-05.03 16:00"
$ws.Range("E37").Value = "2024.05.01 10:00-05.02 17:00"
$ws.Range("F37").Value = 99
$ws.Range("E38").Value = "2024.05.01 10:30-05.02 17:00"
$ws.Range("F38").Value = 1266
$ws.Range("E39").Value = "2024.05.01 19:30-05.01 21:00"
$ws.Range("E40").Value = "2024.05.02 10:00-05.02 17:00"
$ws.Range("F40").Value = 99
$ws.Range("E41").Value = "2024.05.02 10:00-05.03 17:00"
$ws.Range("F41").Value = 1448
$ws.Range("E42").Value = "2024.05.02 10:00-05.02 17:00"
$ws.Range("F42").Value = 18
$ws.Range("E43").Value = "2024.05.12 19:30-05.12 21:30"
$ws.Range("E44").Value = "2024.05.18 10:00-05.18 21:00"
$ws.Range("F44").Value = 328
$ws.Range("E45").Value = "2024.05.25 10:00-05.25 17:00"
$ws.Range("F45").Value = 18
$ws.Range("E46").Value = "2024.06.08 10:00-06.10 17:00"
$ws.Range("F46").Value = 199
$ws.Range("E47").Value = "2024.06.09 10:00-06.09 23:00"
$ws.Range("F47").Value = 706
$ws.Range("E48").Value = "2024.07.20 13:00-07.20 19:00"
$ws.Range("F48").Value = 310
$ws.Range("E49").Value = "2024.11.01 19:30-11.03 17:00"
